{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"],\n  [\"189\u00d79=1701\", \"152\u00d74=608\"],\n  [\"267\u00d74=1068\", \"288\u00d79=2592\"],\n  [\"761\u00d76=4566\", \"220\u00d72=440\"],\n  [\"973\u00d76=5838\", \"249\u00d73=747\"],\n  [\"283\u00d75=1415\", \"387\u00d72=774\"],\n  [\"991\u00d79=8919\", \"470\u00d76=2820\"],\n  [\"958\u00d76=5748\", \"851\u00d77=5957\"],\n  [\"698\u00d74=2792\", \"702\u00d75=3510\"],\n  [\"901\u00d72=1802\", \"654\u00d75=3270\"],\n  [\"713\u00d72=1426\", \"647\u00d75=3235\"],\n  [\"911\u00d78=7288\", \"379\u00d73=1137\"],\n  [\"729\u00d76=4374\", \"106\u00d77=742\"],\n  [\"131\u00d72=262\", \"641\u00d73=1923\"],\n  [\"231\u00d76=1386\", \"120\u00d76=720\"],\n  [\"769\u00d73=2307\", \"370\u00d77=2590\"],\n  [\"463\u00d74=1852\", \"448\u00d78=3584\"],\n  [\"598\u00d75=2990\", \"541\u00d77=3787\"],\n  [\"969\u00d78=7752\", \"532\u00d77=3724\"],\n  [\"444\u00d79=3996\", \"393\u00d73=1179\"],\n  [\"910\u00d75=4550\", \"947\u00d78=7576\"],\n  [\"492\u00d77=3444\", \"406\u00d73=1218\"],\n  [\"429\u00d79=3861\", \"157\u00d72=314\"],\n  [\"396\u00d73=1188\", \"704\u00d76=4224\"],\n  [\"555\u00d74=2220\", \"581\u00d79=5229\"],\n  [\"295\u00d75=1475\", \"926\u00d76=5556\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"),\n    @(\"189\u00d79=1701\", \"152\u00d74=608\"),\n    @(\"267\u00d74=1068\", \"288\u00d79=2592\"),\n    @(\"761\u00d76=4566\", \"220\u00d72=440\"),\n    @(\"973\u00d76=5838\", \"249\u00d73=747\"),\n    @(\"283\u00d75=1415\", \"387\u00d72=774\"),\n    @(\"991\u00d79=8919\", \"470\u00d76=2820\"),\n    @(\"958\u00d76=5748\", \"851\u00d77=5957\"),\n    @(\"698\u00d74=2792\", \"702\u00d75=3510\"),\n    @(\"901\u00d72=1802\", \"654\u00d75=3270\"),\n    @(\"713\u00d72=1426\", \"647\u00d75=3235\"),\n    @(\"911\u00d78=7288\", \"379\u00d73=1137\"),\n    @(\"729\u00d76=4374\", \"106\u00d77=742\"),\n    @(\"131\u00d72=262\", \"641\u00d73=1923\"),\n    @(\"231\u00d76=1386\", \"120\u00d76=720\"),\n    @(\"769\u00d73=2307\", \"370\u00d77=2590\"),\n    @(\"463\u00d74=1852\", \"448\u00d78=3584\"),\n    @(\"598\u00d75=2990\", \"541\u00d77=3787\"),\n    @(\"969\u00d78=7752\", \"532\u00d77=3724\"),\n    @(\"444\u00d79=3996\", \"393\u00d73=1179\"),\n    @(\"910\u00d75=4550\", \"947\u00d78=7576\"),\n    @(\"492\u00d77=3444\", \"406\u00d73=1218\"),\n    @(\"429\u00d79=3861\", \"157\u00d72=314\"),\n    @(\"396\u00d73=1188\", \"704\u00d76=4224\"),\n    @(\"555\u00d74=2220\", \"581\u00d79=5229\"),\n    @(\"295\u00d75=1475\", \"926\u00d76=5556\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $null, 2)\n}\n"}
